# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-39, replacing the prior Strike# values.
$newK = @{
    2  = 6
    3  = 7
    4  = 7
    5  = 6
    6  = 8
    7  = 3
    8  = 3
    9  = 4
    10 = 4
    11 = 3
    12 = 5
    13 = 5
    14 = 2
    15 = 4
    16 = 4
    17 = 7
    18 = 2
    19 = 7
    20 = 9
    21 = 8
    22 = 5
    23 = 4
    24 = 4
    25 = 5
    26 = 6
    27 = 6
    28 = 6
    29 = 4
    30 = 6
    31 = 6
    32 = 6
    33 = 5
    34 = 7
    35 = 5
    36 = 6
    37 = 2
    38 = 4
    39 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
